$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "80-28="
$t.Cell(1,2).Range.Text = "20+44="
$t.Cell(1,3).Range.Text = "51-20="
$t.Cell(1,4).Range.Text = "72-29="
$t.Cell(1,5).Range.Text = "48+1="
$t.Cell(2,1).Range.Text = "11+23="
$t.Cell(2,2).Range.Text = "49-46="
$t.Cell(2,3).Range.Text = "32-18="
$t.Cell(2,4).Range.Text = "14+68="
$t.Cell(2,5).Range.Text = "96-53="
$t.Cell(3,1).Range.Text = "37-16="
$t.Cell(3,2).Range.Text = "43+51="
$t.Cell(3,3).Range.Text = "52+9="
$t.Cell(3,4).Range.Text = "54-16="
$t.Cell(3,5).Range.Text = "62-18="
$t.Cell(4,1).Range.Text = "46+32="
$t.Cell(4,2).Range.Text = "36+13="
$t.Cell(4,3).Range.Text = "62+30="
$t.Cell(4,4).Range.Text = "70-69="
$t.Cell(4,5).Range.Text = "57-0="
$t.Cell(5,1).Range.Text = "71+20="
$t.Cell(5,2).Range.Text = "17+4="
$t.Cell(5,3).Range.Text = "83-72="
$t.Cell(5,4).Range.Text = "84-41="
$t.Cell(5,5).Range.Text = "31+27="
$t.Cell(6,1).Range.Text = "70-41="
$t.Cell(6,2).Range.Text = "89+10="
$t.Cell(6,3).Range.Text = "96-51="
$t.Cell(6,4).Range.Text = "94-54="
$t.Cell(6,5).Range.Text = "19+70="
$t.Cell(7,1).Range.Text = "85-63="
$t.Cell(7,2).Range.Text = "69+13="
$t.Cell(7,3).Range.Text = "96-57="
$t.Cell(7,4).Range.Text = "24+58="
$t.Cell(7,5).Range.Text = "34+38="
$t.Cell(8,1).Range.Text = "23+63="
$t.Cell(8,2).Range.Text = "34+0="
$t.Cell(8,3).Range.Text = "10+66="
$t.Cell(8,4).Range.Text = "3+47="
$t.Cell(8,5).Range.Text = "29+23="
$t.Cell(9,1).Range.Text = "31+14="
$t.Cell(9,2).Range.Text = "40-27="
$t.Cell(9,3).Range.Text = "74+15="
$t.Cell(9,4).Range.Text = "68+17="
$t.Cell(9,5).Range.Text = "39+18="
$t.Cell(10,1).Range.Text = "22+27="
$t.Cell(10,2).Range.Text = "21+72="
$t.Cell(10,3).Range.Text = "39+2="
$t.Cell(10,4).Range.Text = "60-11="
$t.Cell(10,5).Range.Text = "44+23="
$t.Cell(11,1).Range.Text = "4+33="
$t.Cell(11,2).Range.Text = "62+22="
$t.Cell(11,3).Range.Text = "39-24="
$t.Cell(11,4).Range.Text = "35+52="
$t.Cell(11,5).Range.Text = "42+8="
$t.Cell(12,1).Range.Text = "82+6="
$t.Cell(12,2).Range.Text = "53+4="
$t.Cell(12,3).Range.Text = "20+28="
$t.Cell(12,4).Range.Text = "85-51="
$t.Cell(12,5).Range.Text = "32+48="
$t.Cell(13,1).Range.Text = "62+17="
$t.Cell(13,2).Range.Text = "99-91="
$t.Cell(13,3).Range.Text = "44+9="
$t.Cell(13,4).Range.Text = "31+61="
$t.Cell(13,5).Range.Text = "98-51="
$t.Cell(14,1).Range.Text = "95-68="
$t.Cell(14,2).Range.Text = "72-13="
$t.Cell(14,3).Range.Text = "88-65="
$t.Cell(14,4).Range.Text = "67-39="
$t.Cell(14,5).Range.Text = "51+19="
$t.Cell(15,1).Range.Text = "34+10="
$t.Cell(15,2).Range.Text = "39+15="
$t.Cell(15,3).Range.Text = "28+20="
$t.Cell(15,4).Range.Text = "12+65="
$t.Cell(15,5).Range.Text = "67+10="
$t.Cell(16,1).Range.Text = "21+45="
$t.Cell(16,2).Range.Text = "87+5="
$t.Cell(16,3).Range.Text = "74-27="
$t.Cell(16,4).Range.Text = "27-3="
$t.Cell(16,5).Range.Text = "50-18="
$t.Cell(17,1).Range.Text = "26+29="
$t.Cell(17,2).Range.Text = "26+71="
$t.Cell(17,3).Range.Text = "95-49="
$t.Cell(17,4).Range.Text = "4+26="
$t.Cell(17,5).Range.Text = "94-53="
$t.Cell(18,1).Range.Text = "16+49="
$t.Cell(18,2).Range.Text = "47+36="
$t.Cell(18,3).Range.Text = "76-48="
$t.Cell(18,4).Range.Text = "80-52="
$t.Cell(18,5).Range.Text = "63-39="
$t.Cell(19,1).Range.Text = "99-54="
$t.Cell(19,2).Range.Text = "64-63="
$t.Cell(19,3).Range.Text = "9+10="
$t.Cell(19,4).Range.Text = "36+63="
$t.Cell(19,5).Range.Text = "58-14="
$t.Cell(20,1).Range.Text = "53+16="
$t.Cell(20,2).Range.Text = "75-4="
$t.Cell(20,3).Range.Text = "81+3="
$t.Cell(20,4).Range.Text = "12+16="
$t.Cell(20,5).Range.Text = "12+34="
